# Realestate Update resale numbers 2024-01-23 22:23
# Appends a new data row (row 89) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

# Columns A-D hold text that otherwise looks like a date/time/number
# ("2024-01-23", "22:23:32", "03"). Force text interpretation so Excel
# doesn't silently convert them into date/time serials or strip the
# leading zero, then drop the temporary number format again so the new
# row ends up styled the same (i.e. unstyled) as the rest of the table.
$ws.Range("A$row`:D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-23"
$ws.Cells.Item($row, 2).Value = "22:23:32"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "03"

$ws.Range("A$row`:D$row").ClearFormats()

# Columns E-T are plain numeric resale counts (-1 denotes "no data").
$ws.Cells.Item($row, 5).Value = 138480
$ws.Cells.Item($row, 6).Value = 141294
$ws.Cells.Item($row, 7).Value = 171141
$ws.Cells.Item($row, 8).Value = 148877
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 123390
$ws.Cells.Item($row, 11).Value = 223731
$ws.Cells.Item($row, 12).Value = 256274
$ws.Cells.Item($row, 13).Value = 185040
$ws.Cells.Item($row, 14).Value = 110222
$ws.Cells.Item($row, 15).Value = 41342
$ws.Cells.Item($row, 16).Value = 30894
$ws.Cells.Item($row, 17).Value = 73463
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42648
$ws.Cells.Item($row, 20).Value = -1
